$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 39, shifting rows 39..64 down to 40..65.
$ws.Rows.Item(39).Insert(-4121)

# Copy the date-cell formatting from the row above (row 38, which is now above the
# new blank row 39) onto the new row's date cell only, so column D keeps its date style.
$ws.Cells.Item(38, 4).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4122)

# Populate the new row 39 with its values.
$ws.Cells.Item(39, 1).Value = 11
$ws.Cells.Item(39, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(39, 3).Value = "Bíobío"
$ws.Cells.Item(39, 4).Value = 44651
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = 100112012
$ws.Cells.Item(39, 7).Value = "Espinaca"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 250
$ws.Cells.Item(39, 11).Value = 6000
$ws.Cells.Item(39, 12).Value = 6500
$ws.Cells.Item(39, 13).Value = 6200
$ws.Cells.Item(39, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(39, 15).Value = "Región Metropolitana"
$ws.Cells.Item(39, 16).Value = 620
$ws.Cells.Item(39, 17).Value = 10
$ws.Cells.Item(39, 18).Value = "Hortaliza"
